$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 45, shifting existing rows 45-146 down to 46-147
$ws.Rows.Item(45).Insert()

# Populate the new row 45 with the new record data
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44775
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112043
$ws.Cells.Item(45, 7).Value = "Pepino ensalada"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 200
$ws.Cells.Item(45, 11).Value = 18000
$ws.Cells.Item(45, 12).Value = 19000
$ws.Cells.Item(45, 13).Value = 18500
$ws.Cells.Item(45, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 308
$ws.Cells.Item(45, 17).Value = 60
$ws.Cells.Item(45, 18).Value = "Hortaliza"
